# TC_53550 - "Updated base change test case data"
#
# Changes applied:
#  1. Header I7: "Assign Base/Default Base" -> "Assign Base/Default Base Row"
#  2. I8/I9 change from the free-text "5BI 5" [517.050.018] & 801RIL" label
#     to numeric Base Property Index values (14 and 11 respectively), with
#     that text preserved as a cell comment instead.
#  3. J8/J9 Base Property Index values change from 13 to 6.
#  4. Selection / active cell moves from F1:F2 to H9 (and the view scrolls
#     so column C is the left-most visible column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column header text change -----------------------------------------
$ws.Range("I7").Value = "Assign Base/Default Base Row"

# --- Data changes on rows 8 & 9 -----------------------------------------
$ws.Range("I8").Value = 14
$ws.Range("J8").Value = 6

$ws.Range("I9").Value = 11
$ws.Range("J9").Value = 6

# --- Re-attach the old descriptive text as cell comments ----------------
$commentText = "Alpesh Dhakad:" + [char]10 + "5BI 5`" [517.050.018] & 801RIL"

$ws.Range("I8").AddComment($commentText)
$ws.Range("I9").AddComment($commentText)

# --- View / selection state ---------------------------------------------
$excel.ActiveWindow.ScrollColumn = 3
$null = $ws.Range("H9").Select()
